# Update LR-pair table with new TPM-derived values.
# The sending-cluster set grew from {ECs, FAPs, Inflammatory-Mac, Resolving-Mac}
# to {FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac} x {ECs, Inflammatory-Mac, Resolving-Mac}
# i.e. 12 data rows (was 9), all numeric columns recomputed with the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T, rows 2:13 (row 1 header is unchanged)
$data = New-Object 'object[,]' 12,20
# row 2
$data[0,0] = "FAPs"
$data[0,1] = "Tnfsf8"
$data[0,2] = "Tnfrsf8"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.7870423333333333
$data[0,7] = 2.361127
$data[0,8] = 0.2114922334752252
$data[0,9] = 0.2114922334752252
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.7172796666666668
$data[0,13] = 2.151839
$data[0,14] = 0.6605750256943618
$data[0,15] = 0.6605750256943617
$data[0,16] = 0.564529462505889
$data[0,17] = 5.080765162553
$data[0,18] = 0.1397064875620549
$data[0,19] = 0.1397064875620549
# row 3
$data[1,0] = "FAPs"
$data[1,1] = "Tnfsf8"
$data[1,2] = "Tnfrsf8"
$data[1,3] = "Inflammatory-Mac"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.7870423333333333
$data[1,7] = 2.361127
$data[1,8] = 0.2114922334752252
$data[1,9] = 0.2114922334752252
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.2168863333333333
$data[1,13] = 0.650659
$data[1,14] = 0.1997403549444302
$data[1,15] = 0.1997403549444302
$data[1,16] = 0.1706987258547777
$data[1,17] = 1.536288532693
$data[1,18] = 0.04224353378233178
$data[1,19] = 0.04224353378233178
# row 4
$data[2,0] = "FAPs"
$data[2,1] = "Tnfsf8"
$data[2,2] = "Tnfrsf8"
$data[2,3] = "Resolving-Mac"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.7870423333333333
$data[2,7] = 2.361127
$data[2,8] = 0.2114922334752252
$data[2,9] = 0.2114922334752252
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.1516753333333334
$data[2,13] = 0.455026
$data[2,14] = 0.1396846193612081
$data[2,15] = 0.1396846193612081
$data[2,16] = 0.1193749082557778
$data[2,17] = 1.074374174302
$data[2,18] = 0.02954221213083859
$data[2,19] = 0.02954221213083859
# row 5
$data[3,0] = "Inflammatory-Mac"
$data[3,1] = "Tnfsf8"
$data[3,2] = "Tnfrsf8"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 2.080849333333334
$data[3,7] = 6.242548
$data[3,8] = 0.5591611205565395
$data[3,9] = 0.5591611205565395
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.7172796666666668
$data[3,13] = 2.151839
$data[3,14] = 0.6605750256943618
$data[3,15] = 0.6605750256943617
$data[3,16] = 1.492550916196889
$data[3,17] = 13.432958245772
$data[3,18] = 0.3693678715789243
$data[3,19] = 0.3693678715789241
# row 6
$data[4,0] = "Inflammatory-Mac"
$data[4,1] = "Tnfsf8"
$data[4,2] = "Tnfrsf8"
$data[4,3] = "Inflammatory-Mac"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 2.080849333333334
$data[4,7] = 6.242548
$data[4,8] = 0.5591611205565395
$data[4,9] = 0.5591611205565395
$data[4,10] = 1
$data[4,11] = 0.3333333333333333
$data[4,12] = 0.2168863333333333
$data[4,13] = 0.650659
$data[4,14] = 0.1997403549444302
$data[4,15] = 0.1997403549444302
$data[4,16] = 0.4513077821257778
$data[4,17] = 4.061770039132
$data[4,18] = 0.1116870406910885
$data[4,19] = 0.1116870406910885
# row 7
$data[5,0] = "Inflammatory-Mac"
$data[5,1] = "Tnfsf8"
$data[5,2] = "Tnfrsf8"
$data[5,3] = "Resolving-Mac"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 2.080849333333334
$data[5,7] = 6.242548
$data[5,8] = 0.5591611205565395
$data[5,9] = 0.5591611205565395
$data[5,10] = 2
$data[5,11] = 0.6666666666666666
$data[5,12] = 0.1516753333333334
$data[5,13] = 0.455026
$data[5,14] = 0.1396846193612081
$data[5,15] = 0.1396846193612081
$data[5,16] = 0.3156135162497778
$data[5,17] = 2.840521646248
$data[5,18] = 0.07810620828652683
$data[5,19] = 0.07810620828652681
# row 8
$data[6,0] = "MuSCs"
$data[6,1] = "Tnfsf8"
$data[6,2] = "Tnfrsf8"
$data[6,3] = "ECs"
$data[6,4] = 1
$data[6,5] = 0.3333333333333333
$data[6,6] = 0.01097433333333333
$data[6,7] = 0.032923
$data[6,8] = 0.002948998000829621
$data[6,9] = 0.002948998000829622
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.7172796666666668
$data[6,13] = 2.151839
$data[6,14] = 0.6605750256943618
$data[6,15] = 0.6605750256943617
$data[6,16] = 0.007871666155222223
$data[6,17] = 0.07084499539700001
$data[6,18] = 0.001948034430170648
$data[6,19] = 0.001948034430170648
# row 9
$data[7,0] = "MuSCs"
$data[7,1] = "Tnfsf8"
$data[7,2] = "Tnfrsf8"
$data[7,3] = "Inflammatory-Mac"
$data[7,4] = 1
$data[7,5] = 0.3333333333333333
$data[7,6] = 0.01097433333333333
$data[7,7] = 0.032923
$data[7,8] = 0.002948998000829621
$data[7,9] = 0.002948998000829622
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.2168863333333333
$data[7,13] = 0.650659
$data[7,14] = 0.1997403549444302
$data[7,15] = 0.1997403549444302
$data[7,16] = 0.002380182917444444
$data[7,17] = 0.021421646257
$data[7,18] = 0.0005890339074161235
$data[7,19] = 0.0005890339074161236
# row 10
$data[8,0] = "MuSCs"
$data[8,1] = "Tnfsf8"
$data[8,2] = "Tnfrsf8"
$data[8,3] = "Resolving-Mac"
$data[8,4] = 1
$data[8,5] = 0.3333333333333333
$data[8,6] = 0.01097433333333333
$data[8,7] = 0.032923
$data[8,8] = 0.002948998000829621
$data[8,9] = 0.002948998000829622
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.1516753333333334
$data[8,13] = 0.455026
$data[8,14] = 0.1396846193612081
$data[8,15] = 0.1396846193612081
$data[8,16] = 0.001664535666444445
$data[8,17] = 0.014980820998
$data[8,18] = 0.0004119296632428493
$data[8,19] = 0.0004119296632428493
# row 11
$data[9,0] = "Resolving-Mac"
$data[9,1] = "Tnfsf8"
$data[9,2] = "Tnfrsf8"
$data[9,3] = "ECs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.842511
$data[9,7] = 2.527533
$data[9,8] = 0.2263976479674056
$data[9,9] = 0.2263976479674056
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 0.7172796666666668
$data[9,13] = 2.151839
$data[9,14] = 0.6605750256943618
$data[9,15] = 0.6605750256943617
$data[9,16] = 0.6043160092430001
$data[9,17] = 5.438844083187001
$data[9,18] = 0.149552632123212
$data[9,19] = 0.149552632123212
# row 12
$data[10,0] = "Resolving-Mac"
$data[10,1] = "Tnfsf8"
$data[10,2] = "Tnfrsf8"
$data[10,3] = "Inflammatory-Mac"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 0.842511
$data[10,7] = 2.527533
$data[10,8] = 0.2263976479674056
$data[10,9] = 0.2263976479674056
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.2168863333333333
$data[10,13] = 0.650659
$data[10,14] = 0.1997403549444302
$data[10,15] = 0.1997403549444302
$data[10,16] = 0.182729121583
$data[10,17] = 1.644562094247
$data[10,18] = 0.04522074656359373
$data[10,19] = 0.04522074656359375
# row 13
$data[11,0] = "Resolving-Mac"
$data[11,1] = "Tnfsf8"
$data[11,2] = "Tnfrsf8"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 0.842511
$data[11,7] = 2.527533
$data[11,8] = 0.2263976479674056
$data[11,9] = 0.2263976479674056
$data[11,10] = 2
$data[11,11] = 0.6666666666666666
$data[11,12] = 0.1516753333333334
$data[11,13] = 0.455026
$data[11,14] = 0.1396846193612081
$data[11,15] = 0.1396846193612081
$data[11,16] = 0.127788136762
$data[11,17] = 1.150093230858
$data[11,18] = 0.03162426928059984
$data[11,19] = 0.03162426928059984

$ws.Range("A2:T13").Value = $data

